# Replace the old localhost dev URLs with the production domain
# (t-h-logistics.com) across every cell in the "Data" sheet.
#
# All affected values live in the URL column (column B) of the sheet,
# but we scan the whole UsedRange so the edit is robust to layout and
# still only touches cells that actually contain the old host.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHost = "http://localhost:3000"
$newHost = "https://t-h-logistics.com"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$replaced = 0

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldHost)) {
            $cell.Value = $val.Replace($oldHost, $newHost)
            $replaced = $replaced + 1
        }
    }
}

Write-Output "Updated $replaced cell(s) from $oldHost to $newHost"
